# Updates the Coin / Link / Price / Volume(1h) columns with the latest
# coinranking.com snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether it must be forced to
# text (Excel would otherwise read a plain decimal-looking price as a number).
$updates = @(
    @{ Cell = 'D2'; Value = '30.803.93'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +1.95%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.893.20'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +1.24%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '244.98'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +4.37%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4788'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  +1.81%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2903'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +2.05%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '42.93'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +2.99%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.06568'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +0.18%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '21.33'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.42%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.07780'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.40%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.909.15'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  +2.34%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '97.06'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.7408'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +7.32%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '5.182'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +1.91%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '281.38'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +5.56%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '30.789.52'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +1.92%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '13.52'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -1.56%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000007596'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.56%  '; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '2.159.67'; ForceText = $false }
    @{ Cell = 'E22'; Value = '  +2.94%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '5.308'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +1.33%  '; ForceText = $false }
    @{ Cell = 'E24'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '6.244'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +1.46%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '9.366'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -1.52%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '166.38'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +0.40%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '19.14'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  +2.18%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '1.981'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +2.50%  '; ForceText = $false }
    @{ Cell = 'B30'; Value = 'Toncoin'; ForceText = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.374'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +0.23%  '; ForceText = $false }
    @{ Cell = 'B31'; Value = 'Stellar'; ForceText = $false }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    @{ Cell = 'D31'; Value = '0.1001'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +0.91%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '1.518'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +4.37%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '4.377'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +0.70%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '4.122'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +1.93%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.04781'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +0.95%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '1.132'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.25%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.7048'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +0.74%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '2.718'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.01875'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.72%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -0.38%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '6.430'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +2.85%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '70.81'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -2.60%  '; ForceText = $false }
    @{ Cell = 'B43'; Value = 'TheSandbox'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.4212'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +1.42%  '; ForceText = $false }
    @{ Cell = 'B44'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = 'D44'; Value = '1.926'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.8486'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.00%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '102.52'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -0.40%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '9.436'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +3.72%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  +1.33%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '942.84'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -3.41%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '35.29'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +2.43%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.ForceText) {
        # Leading apostrophe = Excel's own 'treat as text' quote prefix,
        # so numeric-looking strings (e.g. "244.98") are kept as text,
        # exactly like the un-formatted price/volume columns in this sheet.
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
